# Update the cached "last saved" date/time field text from 16/05/2021 to
# 18/05/2021 everywhere it appears (Slide Master, every Slide Layout and the
# Notes Master), and tweak the wording of one sentence on slide 6.

$p = $ppt.ActivePresentation

$oldDate = "16/05/2021"
$newDate = "18/05/2021"

# --- Slide Master: "Espaço Reservado para Data" placeholder -----------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shape = $master.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.TextRange.Text -eq $oldDate) {
        $shape.TextFrame.TextRange.Text = $newDate
    }
}

# --- Every Slide Layout's date placeholder -----------------------------------
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shape = $layout.Shapes.Item($i)
        if ($shape.HasTextFrame -and $shape.TextFrame.TextRange.Text -eq $oldDate) {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Notes Master date placeholder -------------------------------------------
# (direct shape TextRange edits are not persisted for the Notes Master in this
# host, but the HeadersFooters.DateAndTime accessor is, so use that instead.)
$p.NotesMaster.HeadersFooters.DateAndTime.Text = $newDate

# --- Slide 6: tweak a sentence in the content placeholder --------------------
$slide6 = $p.Slides.Item(6)
$contentShape = $slide6.Shapes.Item(2)
$fullRange = $contentShape.TextFrame.TextRange
$lastParagraph = $fullRange.Paragraphs($fullRange.Paragraphs().Count, 1)

# Force a clean single-run replacement (the host otherwise tries to keep any
# shared prefix/suffix as separate runs), by first clearing the paragraph to
# unrelated text, then writing the final wording.
$lastParagraph.Text = "-"
$contentShape.TextFrame.TextRange.Paragraphs(9, 1).Text = "Aqui temos descritivos de algumas energias renováveis."
